$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The row for rule "R30" (B10) had its "Integer min" (column C) value
# restored from 18 back to 1.
$ws.Range("C10").Value = 1
